# Generate Report for Handback
# Updates the localization-status workbook with the results for the
# 682f6a55-dfd7-4031-9ea0-3c708bb948b6 handback file (row 7) on both the
# zh-cn and de-de sheets: a Latest Target File hyperlink, a Latest
# Handback File name, a Latest Handback DateTime, and an Error Detail
# message describing that the handback was not done against the latest
# version. Also widens the Error Detail column.

$wb = $excel.ActiveWorkbook

$mdDisplay = "682f6a55-dfd7-4031-9ea0-3c708bb948b6.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2068644d2ace7869b2eb5f3d3de58f49e166c2/e2e/682f6a55-dfd7-4031-9ea0-3c708bb948b6.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c78ac907fb2d1e945d8bd17f8f21b784d9e0817c/e2e/682f6a55-dfd7-4031-9ea0-3c708bb948b6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a2068644d2ace7869b2eb5f3d3de58f49e166c2/e2e/682f6a55-dfd7-4031-9ea0-3c708bb948b6.md."

function Update-HandbackRow {
    param(
        $ws,
        [string]$handbackFile,
        [string]$handbackDateTime
    )

    # Latest Target File: hyperlink to the handback markdown file on GitHub
    $ws.Hyperlinks.Add($ws.Range("I7"), $mdUrl, "", "", $mdDisplay)

    # Latest Handback File
    $ws.Range("J7").Value = $handbackFile

    # Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime

    # Error Detail
    $ws.Range("P7").Value = $errorDetail

    # Error Detail column is now much wider to fit the message
    $ws.Range("P1").ColumnWidth = 39.17
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow -ws $wsZhCn `
    -handbackFile "682f6a55-dfd7-4031-9ea0-3c708bb948b6.e0d6d0db005523d8da2b4f7574e96bd5d9ae0093.zh-cn.xlf" `
    -handbackDateTime "2016-09-02 08:49:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow -ws $wsDeDe `
    -handbackFile "682f6a55-dfd7-4031-9ea0-3c708bb948b6.e0d6d0db005523d8da2b4f7574e96bd5d9ae0093.de-de.xlf" `
    -handbackDateTime "2016-09-02 08:49:52"
